# Auto-generated edit script applying the diff to Seraph_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("H12").Value = 201
$ws.Range("I12").Value = 200
$ws.Range("K12").Value = 200
$ws.Range("M12").Value = -30
$ws.Range("H15").Value = 806.13336
$ws.Range("I15").Value = 806.13336
$ws.Range("K15").Value = 2418.40008
$ws.Range("M15").Value = -2249.40008
$ws.Range("H55").Value = 1632.8889
$ws.Range("J55").Value = 2966
$ws.Range("L55").Value = 2966
$ws.Range("N55").Value = -3394

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3013.6
$ws.Range("I2").Value = 1757
$ws.Range("J2").Value = 4898.5
$ws.Range("K2").Value = 1757
$ws.Range("L2").Value = 4898.5
$ws.Range("M2").Value = -1644
$ws.Range("N2").Value = -5124.5
$ws.Range("H3").Value = 10006
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("H5").Value = 244.22223
$ws.Range("I5").Value = 57.333332
$ws.Range("J5").Value = 337.66666
$ws.Range("K5").Value = 57.333332
$ws.Range("L5").Value = 337.66666
$ws.Range("M5").Value = 54.666668
$ws.Range("N5").Value = -561.66666
$ws.Range("H32").Value = 4369.104
$ws.Range("I32").Value = 3152.027
$ws.Range("K32").Value = 3152.027
$ws.Range("M32").Value = -2865.027
$ws.Range("H35").Value = 26341.75
$ws.Range("I35").Value = 26341.75
$ws.Range("K35").Value = 26341.75
$ws.Range("M35").Value = -25935.75
$ws.Range("H63").Value = 2605
$ws.Range("I63").Value = 1675
$ws.Range("K63").Value = 1675
$ws.Range("M63").Value = -989
$ws.Range("H66").Value = 2605
$ws.Range("I66").Value = 1675
$ws.Range("K66").Value = 8375
$ws.Range("M66").Value = -4943
$ws.Range("H116").Value = 3013.6
$ws.Range("I116").Value = 1757
$ws.Range("J116").Value = 4898.5
$ws.Range("K116").Value = 1757
$ws.Range("L116").Value = 4898.5
$ws.Range("M116").Value = 537
$ws.Range("N116").Value = -9486.5
$ws.Range("H132").Value = 1903
$ws.Range("I132").Value = 1903
$ws.Range("K132").Value = 5709
$ws.Range("M132").Value = -3179

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3013.6
$ws.Range("I3").Value = 1757
$ws.Range("J3").Value = 4898.5
$ws.Range("K3").Value = 1757
$ws.Range("L3").Value = 4898.5
$ws.Range("M3").Value = -1643
$ws.Range("N3").Value = -5126.5
$ws.Range("H4").Value = 244.22223
$ws.Range("I4").Value = 57.333332
$ws.Range("J4").Value = 337.66666
$ws.Range("K4").Value = 57.333332
$ws.Range("L4").Value = 337.66666
$ws.Range("M4").Value = 57.666668
$ws.Range("N4").Value = -567.66666
$ws.Range("H26").Value = 17900.5
$ws.Range("I26").Value = 17900.5
$ws.Range("K26").Value = 17900.5
$ws.Range("M26").Value = -17608.5
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2776.5
$ws.Range("I31").Value = 3662
$ws.Range("K31").Value = 3662
$ws.Range("M31").Value = -3367
$ws.Range("H34").Value = 2776.5
$ws.Range("I34").Value = 3662
$ws.Range("K34").Value = 3662
$ws.Range("M34").Value = -3460
$ws.Range("H132").Value = 3988.65
$ws.Range("I132").Value = 2611.375
$ws.Range("K132").Value = 7834.125
$ws.Range("M132").Value = -5304.125
$ws.Range("H134").Value = 1808.4445
$ws.Range("I134").Value = 1544.0278
$ws.Range("J134").Value = 2866.111
$ws.Range("K134").Value = 4632.0834
$ws.Range("L134").Value = 8598.332999999999
$ws.Range("M134").Value = -2097.0834
$ws.Range("N134").Value = -13668.333
$ws.Range("H141").Value = 146159.2
$ws.Range("J141").Value = 146159.2
$ws.Range("L141").Value = 146159.2
$ws.Range("N141").Value = -156519.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 341.78946
$ws.Range("J12").Value = 351.61905
$ws.Range("L12").Value = 1054.85715
$ws.Range("N12").Value = -1400.85715
$ws.Range("H126").Value = 2825
$ws.Range("I126").Value = 2825
$ws.Range("K126").Value = 8475
$ws.Range("M126").Value = -3535
$ws.Range("H132").Value = 8441.706
$ws.Range("I132").Value = 9528.615
$ws.Range("K132").Value = 85757.535
$ws.Range("M132").Value = -83227.535

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 1000
$ws.Range("J41").Value = 1000
$ws.Range("L41").Value = 1000
$ws.Range("N41").Value = -1710
$ws.Range("H70").Value = 7870.375
$ws.Range("I70").Value = 7497
$ws.Range("K70").Value = 7497
$ws.Range("M70").Value = -7227
$ws.Range("H73").Value = 7870.375
$ws.Range("I73").Value = 7497
$ws.Range("K73").Value = 7497
$ws.Range("M73").Value = -6561
$ws.Range("H122").Value = 30515.371
$ws.Range("I122").Value = 1905.7084
$ws.Range("J122").Value = 92936.45
$ws.Range("K122").Value = 5717.1252
$ws.Range("L122").Value = 278809.35
$ws.Range("M122").Value = -3267.1252
$ws.Range("N122").Value = -283709.35

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 6434.75
$ws.Range("I22").Value = 5244.533
$ws.Range("K22").Value = 5244.533
$ws.Range("M22").Value = -4949.533
$ws.Range("H27").Value = 6434.75
$ws.Range("I27").Value = 5244.533
$ws.Range("K27").Value = 5244.533
$ws.Range("M27").Value = -5137.533
$ws.Range("H30").Value = 2711.3
$ws.Range("I30").Value = 2735.5
$ws.Range("J30").Value = 2675
$ws.Range("K30").Value = 2735.5
$ws.Range("L30").Value = 2675
$ws.Range("M30").Value = -2627.5
$ws.Range("N30").Value = -2891
$ws.Range("H132").Value = 6468.125
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 6468.125
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 19404.375
$ws.Range("N132").Value = -24464.375
$ws.Range("M132").ClearContents()
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 6923.5
$ws.Range("J52").Value = 4847
$ws.Range("L52").Value = 4847
$ws.Range("N52").Value = -5299
$ws.Range("H136").Value = 1281.3
$ws.Range("I136").Value = 1304.6207
$ws.Range("K136").Value = 3913.8621
$ws.Range("M136").Value = -1363.8621
